$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.392.25'
$ws.Range("E2").Value = '  +2.36%  '
$ws.Range("D3").Value = '1.790.47'
$ws.Range("E3").Value = '  +2.87%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.006'
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '337.30'
$ws.Range("E5").Value = '  +0.99%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.003'
$ws.Range("E6").Value = '  +0.63%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3790'
$ws.Range("E7").Value = '  +1.78%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3456'
$ws.Range("E8").Value = '  +2.19%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '48.67'
$ws.Range("E9").Value = '  +0.79%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.197'
$ws.Range("E10").Value = '  +1.08%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07497'
$ws.Range("E11").Value = '  +0.06%  '
$ws.Range("E12").Value = '  +1.14%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.93'
$ws.Range("E13").Value = '  +7.50%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.458'
$ws.Range("E14").Value = '  +1.53%  '
$ws.Range("D15").Value = '1.791.57'
$ws.Range("E15").Value = '  +2.45%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.062'
$ws.Range("E16").Value = '  +0.42%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001101'
$ws.Range("E17").Value = '  +1.80%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06662'
$ws.Range("E18").Value = '  -0.26%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '84.77'
$ws.Range("E19").Value = '  +2.48%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.002'
$ws.Range("E20").Value = '  +0.48%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.519'
$ws.Range("E21").Value = '  +4.85%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '17.33'
$ws.Range("E22").Value = '  +3.82%  '
$ws.Range("D23").Value = '27.406.28'
$ws.Range("E23").Value = '  +2.14%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.49'
$ws.Range("E24").Value = '  -2.71%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.434'
$ws.Range("E25").Value = '  +0.01%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.557'
$ws.Range("E26").Value = '  +5.91%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.488'
$ws.Range("E27").Value = '  +1.07%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '21.43'
$ws.Range("E28").Value = '  +9.36%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '153.80'
$ws.Range("E29").Value = '  +1.24%  '
$ws.Range("D30").Value = '1.996.44'
$ws.Range("E30").Value = '  +2.78%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '133.49'
$ws.Range("E31").Value = '  +0.99%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.067'
$ws.Range("E32").Value = '  -0.92%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.087'
$ws.Range("E33").Value = '  +1.31%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08683'
$ws.Range("E34").Value = '  +1.33%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '13.18'
$ws.Range("E35").Value = '  +2.33%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.658'
$ws.Range("E36").Value = '  -1.56%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.448'
$ws.Range("E37").Value = '  +0.39%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.6893'
$ws.Range("E38").Value = '  +9.44%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06362'
$ws.Range("E39").Value = '  +1.56%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.845'
$ws.Range("E40").Value = '  +4.02%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.2202'
$ws.Range("E41").Value = '  +1.68%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.02346'
$ws.Range("E42").Value = '  -0.03%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.267'
$ws.Range("E43").Value = '  +4.10%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.39'
$ws.Range("E44").Value = '  +0.37%  '
$ws.Range("B45").Value = 'Frax'
$ws.Range("C45").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.002'
$ws.Range("E45").Value = '  +0.68%  '
$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6416'
$ws.Range("E46").Value = '  +3.24%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.864'
$ws.Range("E47").Value = '  -1.53%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.131'
$ws.Range("E48").Value = '  +2.85%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '129.76'
$ws.Range("E49").Value = '  +0.48%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07203'
$ws.Range("E50").Value = '  -0.30%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '79.43'
$ws.Range("E51").Value = '  +1.71%  '
